$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the row/cell layout below is authoritative
$ws.Cells.Clear()


# Column widths: A and B get new custom widths; D keeps its existing width.
$ws.Columns.Item(1).ColumnWidth = 19.85546875
$ws.Columns.Item(2).ColumnWidth = 18.85546875
$ws.Columns.Item(4).ColumnWidth = 18.28515625

# --- Cell values / formulas ---
$ws.Range("A1").Value = "Pulmonary Edema"
$ws.Range("B1").Value = "Peripheral Edema"
$ws.Range("C1").Value = "# of pts"
$ws.Range("D1").Value = "Percent of total"
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 12705
$ws.Range("D2").Formula = "=C2/C12"
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2011
$ws.Range("D3").Formula = "=C3/C12"
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = -1
$ws.Range("C4").Value = 1512
$ws.Range("D4").Formula = "=C4/C12"
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1904
$ws.Range("D5").Formula = "=C5/C12"
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 838
$ws.Range("D6").Formula = "=C6/C12"
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = -1
$ws.Range("C7").Value = 137
$ws.Range("D7").Formula = "=C7/C12"
$ws.Range("A8").Value = -1
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 762
$ws.Range("D8").Formula = "=C8/C12"
$ws.Range("A9").Value = -1
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 52
$ws.Range("D9").Formula = "=C9/C12"
$ws.Range("A10").Value = -1
$ws.Range("B10").Value = -1
$ws.Range("C10").Value = 724
$ws.Range("D10").Formula = "=C10/C12"
$ws.Range("A11").Value = "No Exam Found"
$ws.Range("C11").Value = 5377
$ws.Range("D11").Formula = "=C11/C12"
$ws.Range("A12").Value = "Total"
$ws.Range("C12").Formula = "=SUM(C2:C11)"
$ws.Range("D12").Formula = "=SUM(D2:D11)"
$ws.Range("C14").Value = "Percent"
$ws.Range("A15").Value = "Missing Info"
$ws.Range("C15").Formula = "=(C4+C7+C8+C9+C10+C11)/C12"
$ws.Range("D15").Value = "Addendums (approx)"
$ws.Range("F15").Value = 2301
$ws.Range("A16").Value = "All info"
$ws.Range("C16").Formula = "=(C2+C3+C5+C6)/C12"
$ws.Range("A17").Value = "Of pts with all info:"
$ws.Range("A18").Value = "Pulmonary Edema"
$ws.Range("B18").Value = "Peripheral Edema"
$ws.Range("D18").Value = "Percent of all info pts"
$ws.Range("A19").Value = 0
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 12705
$ws.Range("D19").Formula = "=C19/C23"
$ws.Range("A20").Value = 0
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 2011
$ws.Range("D20").Formula = "=C20/C23"
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 1904
$ws.Range("D21").Formula = "=C21/C23"
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 838
$ws.Range("D22").Formula = "=C22/C23"
$ws.Range("A23").Value = "Total"
$ws.Range("C23").Formula = "=SUM(C19:C22)"
$ws.Range("D23").Formula = "=SUM(D19:D22)"
$ws.Range("A25").Value = "REVISED: NEONATES REMOVED."
$ws.Range("A26").Value = "Pulmonary Edema"
$ws.Range("B26").Value = "Peripheral Edema"
$ws.Range("C26").Value = "#patients"
$ws.Range("D26").Value = "Percent total"
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 838
$ws.Range("D27").Formula = "=C27/C37"
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 1751
$ws.Range("D28").Formula = "=C28/C37"
$ws.Range("A29").Value = 1
$ws.Range("B29").Value = -1
$ws.Range("C29").Value = 96
$ws.Range("D29").Formula = "=C29/C37"
$ws.Range("A30").Value = 0
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 2009
$ws.Range("D30").Formula = "=C30/C37"
$ws.Range("A31").Value = 0
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 11308
$ws.Range("D31").Formula = "=C31/C37"
$ws.Range("A32").Value = 0
$ws.Range("B32").Value = -1
$ws.Range("C32").Value = 1145
$ws.Range("D32").Formula = "=C32/C37"
$ws.Range("A33").Value = -1
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = 52
$ws.Range("D33").Formula = "=C33/C37"
$ws.Range("A34").Value = -1
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 539
$ws.Range("D34").Formula = "=C34/C37"
$ws.Range("A35").Value = -1
$ws.Range("B35").Value = -1
$ws.Range("C35").Value = 557
$ws.Range("D35").Formula = "=C35/C37"
$ws.Range("A36").Value = "No Exam"
$ws.Range("C36").Value = 4280
$ws.Range("D36").Formula = "=C36/C37"
$ws.Range("A37").Value = "Total"
$ws.Range("C37").Formula = "=SUM(C27:C36)"
$ws.Range("D37").Formula = "=SUM(D27:D36)"
$ws.Range("B38").Value = "Percent"
$ws.Range("A39").Value = "Missing Info"
$ws.Range("B39").Formula = "=(C29+C32+C33+C34+C35+C36)/C37"
$ws.Range("A40").Value = "All Info"
$ws.Range("B40").Formula = "=1-B39"
$ws.Range("A41").Value = "Of pts with all info:"
$ws.Range("A42").Value = "Pulmonary Edema"
$ws.Range("B42").Value = "Peripheral Edema"
$ws.Range("C42").Value = "#"
$ws.Range("D42").Value = "Percent total"
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 838
$ws.Range("D43").Formula = "=C43/C47"
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = 0
$ws.Range("C44").Value = 1751
$ws.Range("D44").Formula = "=C44/C47"
$ws.Range("A45").Value = 0
$ws.Range("B45").Value = 1
$ws.Range("C45").Value = 2009
$ws.Range("D45").Formula = "=C45/C47"
$ws.Range("A46").Value = 0
$ws.Range("B46").Value = 0
$ws.Range("C46").Value = 11308
$ws.Range("D46").Formula = "=C46/C47"
$ws.Range("A47").Value = "Total"
$ws.Range("C47").Formula = "=SUM(C43:C46)"
$ws.Range("D47").Formula = "=SUM(D43:D46)"

# --- Number formats (applied last; a Formula write right after a
#     NumberFormat call on another cell can otherwise inherit the style) ---
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("C3").NumberFormat = "#,##0"
$ws.Range("C4").NumberFormat = "#,##0"
$ws.Range("C5").NumberFormat = "#,##0"
$ws.Range("C7").NumberFormat = "#,##0"
$ws.Range("C8").NumberFormat = "#,##0"
$ws.Range("C9").NumberFormat = "#,##0"
$ws.Range("C10").NumberFormat = "#,##0"
$ws.Range("C11").NumberFormat = "#,##0"
$ws.Range("C12").NumberFormat = "#,##0"
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C23").NumberFormat = "#,##0"

# Selection / scroll position to match the edited view
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D38").Select()
